$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row number -> new dt_insertion serial value (column D)
$values = @{
    2 = 45492.67320601852
    3 = 45492.67325231482
    4 = 45492.67326388889
    5 = 45492.67320601852
    6 = 45492.67322916666
    7 = 45492.67334490741
    8 = 45492.67326388889
    9 = 45492.67321759259
    10 = 45492.67335648148
    11 = 45492.67321759259
    12 = 45492.67335648148
    13 = 45492.67332175926
    14 = 45492.67318287037
    15 = 45492.67324074074
    16 = 45492.67328703704
    17 = 45492.67336805556
    18 = 45492.67329861111
    19 = 45492.67325231482
    20 = 45492.67319444445
    21 = 45492.67328703704
    22 = 45492.67329861111
    23 = 45492.67332175926
    24 = 45492.67331018519
    25 = 45492.67331018519
    26 = 45492.67333333333
    27 = 45492.67333333333
    28 = 45492.67317129629
    29 = 45492.67318287037
    30 = 45492.67327546296
    31 = 45492.67329861111
    32 = 45492.67336805556
    33 = 45492.67334490741
    34 = 45492.67317129629
    35 = 45492.67319444445
    36 = 45492.67322916666
    37 = 45492.67324074074
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 4).Value2 = $values[$row]
}

Write-Output "Updated D column values for rows 2-37"
